$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ------------------------------------------------------------------
# 1) Insert two new rows for the new "2509" period entries, right
#    after the current last data row (22), preserving the table's
#    bottom-border formatting on the new last row.
# ------------------------------------------------------------------
$ws.Range("B23:J24").Insert(-4121)

# carry the special "last row" (thicker bottom border) formatting down to the new last row
$ws.Range("B22:J22").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)   # xlPasteFormats

# row 22 is no longer the last row -> give it the regular middle-row formatting
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats

# new row 23 also gets the regular middle-row formatting
$ws.Range("B21:J21").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Fill in the two new data rows (period 2509) for the same two
#    workers that already had a 2508 entry.
# ------------------------------------------------------------------
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "12598379"
$ws.Range("D23").Value = "ILMER IVAN PASSO PUELLO"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047445983"
$ws.Range("D24").Value = "DAYANA PAOLA ALTAMAR DIAZ"
$ws.Range("E24").Value = "2509"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

# ------------------------------------------------------------------
# 3) Center-align the "Periodo Mora" column for all the data rows.
# ------------------------------------------------------------------
$ws.Range("E16:E23").HorizontalAlignment = -4108   # xlCenter

# ------------------------------------------------------------------
# 4) Update the summary figures: count of periods now 5 (was 4), and
#    total "Valor Mora" now reflects the two additional rows.
# ------------------------------------------------------------------
$ws.Range("F13").Value = 5
$ws.Range("E11").Value = 456731

Write-Host ("Dimension after edit: " + $ws.UsedRange.Address())
